$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New query/label text per commit "CTDC diagnosis 13 scripts"
$statQuery = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
   WHERE c.disease =  "Adenocarcinoma of the gastroesophageal junction"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files'

$filesTabLabel = 'FilesTab'

$casesTabQuery = 'MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
 WHERE c.disease = "Adenocarcinoma of the gastroesophageal junction"
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity'

$filesTabQuery = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
 WHERE c.disease = "Adenocarcinoma of the gastroesophageal junction"
WITH
    f, parent, c, a, ct,
    [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`'

# Row 2 (CasesTab row): swap in the refreshed CasesTab + StatQuery text
$ws.Range("B2").Value = $casesTabQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 195

# Row 3 (new FilesTab row)
$ws.Range("A3").Value = $filesTabLabel
$ws.Range("B3").Value = $filesTabQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 409.5

# Scroll the view down so row 3 is visible, and leave B3 as the active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("B3").Select()
